# Replace every occurrence of "congenital" in column A (across all
# worksheets) with "misc_long_term". The source workbook stores one
# variable-name list per sheet in column A; each of the affected sheets
# has a single "congenital" entry somewhere near the top of the column.
$wb = $excel.ActiveWorkbook

$replacedCount = 0

foreach ($ws in $wb.Worksheets) {
    $used = $ws.UsedRange
    if ($used -eq $null) {
        continue
    }

    $rowCount = $used.Rows.Count
    for ($i = 1; $i -le $rowCount; $i++) {
        $cell = $ws.Cells.Item($i, 1)
        if ($cell.Value2 -eq "congenital") {
            $cell.Value = "misc_long_term"
            $replacedCount = $replacedCount + 1
        }
    }
}

Write-Host "Replaced 'congenital' with 'misc_long_term' in $replacedCount cell(s)."
